# Add a new "Functionality" worksheet as the first sheet in the workbook,
# containing a components/functionality catalog table.

$wb = $excel.ActiveWorkbook

# Insert a brand new worksheet before the first existing sheet ("Ships")
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "Functionality"

# Header row
$newSheet.Range("A1").Value = "Class"
$newSheet.Range("B1").Value = "Details"
$newSheet.Range("C1").Value = "Requirements"
$newSheet.Range("D1").Value = "Data"
$newSheet.Range("E1").Value = "Functionality"

# Data rows
$newSheet.Range("A2").Value = "IDieTarget"
$newSheet.Range("B2").Value = "This is just an optional component if Die() needs to kill a parent."
$newSheet.Range("C2").Value = "-"
$newSheet.Range("D2").Value = "-"
$newSheet.Range("E2").Value = "-"

$newSheet.Range("A3").Value = "Idie"
$newSheet.Range("B3").Value = "Destroys the gameObject, or the IDieTarget gameObject if specified."
$newSheet.Range("C3").Value = "-"
$newSheet.Range("D3").Value = "GameObject that dies"
$newSheet.Range("E3").Value = "Die()"

$newSheet.Range("A4").Value = "Health"
$newSheet.Range("B4").Value = "Tracks the max and current health of an entity."
$newSheet.Range("C4").Value = "-"
$newSheet.Range("D4").Value = "Max, Current"
$newSheet.Range("E4").Value = "Change(amount +/-)"

$newSheet.Range("A5").Value = "ITakeDamage"
$newSheet.Range("B5").Value = "Accepts incoming damage and reduces current health appropriately."
$newSheet.Range("C5").Value = "Health"
$newSheet.Range("D5").Value = "Health, Idie, MinSpeedForImpactDamage, ExtraSpeedDamageMultiplier"
$newSheet.Range("E5").Value = "TakeDamage(amount +, collisionSpeed = 0)"

$newSheet.Range("A6").Value = "IDealImpactDamage"
$newSheet.Range("B6").Value = "Attempts to deal damage to anything that collides with it. Can limit frequency."
$newSheet.Range("C6").Value = "-"
$newSheet.Range("D6").Value = "RigidBody2D, ImpactBaseDamage, DamageInterval"
$newSheet.Range("E6").Value = "OnCollisionEnter2D(dealImpactDamage(ITakeDamage damageableTarget, float collisionSpeed)"

# Formatting: whole sheet uses 8pt Times New Roman, header row is bold
$usedRange = $newSheet.Range("A1:E6")
$usedRange.Font.Name = "Times New Roman"
$usedRange.Font.Size = 8

$headerRange = $newSheet.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.RowHeight = 10.5

# Column widths (target OOXML <col> widths are 20.140625 / 53.42578125 / 12.28515625 /
# 51 / 85.28515625; ColumnWidth is expressed in the "Normal"-style character unit, which
# the exporter re-expresses in XML with a constant ~5/6 character padding, so back it out
# here to land on the intended stored width)
$newSheet.Columns.Item(1).ColumnWidth = 19.307291666666668
$newSheet.Columns.Item(2).ColumnWidth = 52.592447916666664
$newSheet.Columns.Item(3).ColumnWidth = 11.451822916666666
$newSheet.Columns.Item(4).ColumnWidth = 50.166666666666664
$newSheet.Columns.Item(5).ColumnWidth = 84.45182291666667

# Select A8 as the active cell like the target file shows
$newSheet.Range("A8").Select()

# The new sheet should be the active/selected tab
$newSheet.Activate()
